# control de progreso.xlsx - add "estatus" (status) tracking subtasks to item 8,
# mark task 7.07 as completed, and add new task 18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$tbl = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# 1) Mark row 34 (task 7.07 "añadirlo a detallesPedido") as completed ("si")
#    and stamp it with start/end dates, copying the date format used by the
#    neighbouring rows (style s="1").
# ---------------------------------------------------------------------------
$ws.Range("C34").Value2 = "si"

$ws.Range("D33").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E33").Copy() | Out-Null
$ws.Range("E34").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("D34").Value2 = 43803
$ws.Range("E34").Value2 = 43803

# ---------------------------------------------------------------------------
# 2) Remove the old standalone task "9" (row 36: "Agregar estatus a la
#    visualizacion, al lado de fecha de ingreso") - it becomes a sub task of
#    item 8 below instead.
# ---------------------------------------------------------------------------
$ws.Rows.Item(36).Delete()

# ---------------------------------------------------------------------------
# 3) Insert 6 fresh rows right after item "8" (row 35) for the new 8.01-8.06
#    sub tasks.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 6; $i++) {
    $ws.Rows.Item(36).Insert()
}

$subTasks = @(
    @(8.01, "agregarlo a la interfaz", "si"),
    @(8.02, "agregarlo a la base de datos", "si"),
    @(8.03, "que se pueda guardar", "no"),
    @(8.04, "integrarlo a todas las funcionalidades de Pedido", "no"),
    @(8.05, "que se muestre en la visualizacion, a un lado de la fecha de ingreso", "no"),
    @(8.06, "que se pueda modificar en detalles pedido", "no")
)

$r = 36
foreach ($t in $subTasks) {
    $ws.Range("A$r").Value2 = $t[0]
    $ws.Range("B$r").Value2 = $t[1]
    $ws.Range("C$r").Value2 = $t[2]
    $r++
}

# ---------------------------------------------------------------------------
# 4) Append the new task "18" and the new sub task "8.07" at the end of the
#    table (rows 50 and 51).
# ---------------------------------------------------------------------------
$ws.Range("A50").Value2 = 18
$ws.Range("B50").Value2 = "cambiar el formato de la fecha que se muestra en el reporte de produccion"
$ws.Range("C50").Value2 = "no"

$ws.Range("A51").Value2 = 8.07
$ws.Range("B51").Value2 = "quitarlo de partidas"
$ws.Range("C51").Value2 = "no"

# ---------------------------------------------------------------------------
# 5) Resize the table to cover the new data range.
# ---------------------------------------------------------------------------
$tbl.Resize($ws.Range("A1:E51"))

# ---------------------------------------------------------------------------
# 6) Re-apply the "Terminado" (status) filter to show only the pending ("no")
#    rows - this hides every "si" row, same as before the edit.
# ---------------------------------------------------------------------------
$tbl.Range.AutoFilter(3, @("no"), 7) | Out-Null

# ---------------------------------------------------------------------------
# 7) Update the view: drop the frozen scroll position and select C37 instead.
# ---------------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null
$ws.Range("C37").Select() | Out-Null
